# Auto-generated Excel COM-interop edit script
# Refreshes the BRVM "Recommandations" / "Top_YTD" tables to the values
# captured by the latest automated scrape (see commit message:
# "Update automatique BRVM via GitHub Actions").
# Only the cells that actually changed value are touched; row order in the
# sheet (r="n" indices) is preserved exactly as in the source workbook, so
# some rows are rewritten in place to reflect the new sort-by-performance order.

$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------------
# Sheet "Recommandations"
# ----------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Recommandations")

# Row 2: BRVM - SERVICES PUBLICS
$ws1.Range("D2").Value = 83303.74000000001

# Row 3: VIVO ENERGY CI
$ws1.Range("D3").Value = 78420

# Row 4: SUCRIVOIRE
$ws1.Range("C4").Value = 74
$ws1.Range("D4").Value = 69655

# Row 5: BRVM - AUTRES SECTEURS
$ws1.Range("D5").Value = 65217.07

# Row 6: SAFCA CI
$ws1.Range("D6").Value = 61245

# Row 7: NEI-CEDA CI
$ws1.Range("D7").Value = 59760

# Row 8: CFAO MOTORS CI
$ws1.Range("D8").Value = 59295

# Row 10: AIR LIQUIDE CI
$ws1.Range("D10").Value = 49075

# Row 11: UNIWAX CI
$ws1.Range("D11").Value = 43010

# Row 12: BRVM - DISTRIBUTION
$ws1.Range("D12").Value = 36667.86

# Row 13: BRVM - TRANSPORT
$ws1.Range("D13").Value = 32722.81

# Row 14: BRVM - AGRICULTURE
$ws1.Range("D14").Value = 23260.69

# Row 16: BRVM-PRINCIPAL
$ws1.Range("D16").Value = 15205.39

# Row 17: BRVM - INDUSTRIE
$ws1.Range("D17").Value = 13990.53

# Row 18: BRVM - CONSOMMATION DE BASE
$ws1.Range("D18").Value = 12047.67

# Row 19: BRVM-PRESTIGE
$ws1.Range("D19").Value = 11900.99

# Row 20: BRVM - INDUSTRIELS
$ws1.Range("D20").Value = 11469.5

# Row 21: BRVM - ENERGIE
$ws1.Range("D21").Value = 11204.19

# Row 22: BRVM - FINANCES
$ws1.Range("D22").Value = 10998.91

# Row 23: BRVM - SERVICES FINANCIERS
$ws1.Range("D23").Value = 10809.45

# Row 24: BRVM - CONSOMMATION DISCRETIONNAIRE
$ws1.Range("D24").Value = 9684.82

# Row 25: BRVM - TELECOMMUNICATIONS
$ws1.Range("D25").Value = 9448.48

# Row 27: SITAB CI (STBC)
$ws1.Range("B27").Value = 23
$ws1.Range("C27").Value = 7
$ws1.Range("D27").Value = 117.78

# Row 28: FILTISAC CI (FTSC)
$ws1.Range("B28").Value = 22
$ws1.Range("D28").Value = 113.57

# Row 31: UNIWAX CI (UNXC)
$ws1.Range("A31").Value = 'UNIWAX CI (UNXC)'
$ws1.Range("B31").Value = 23
$ws1.Range("C31").Value = 16
$ws1.Range("D31").Value = 50.82
$ws1.Range("E31").Value = 7.37
$ws1.Range("G31").Value = '👀 À surveiller'

# Row 32: ECOBANK TRANS. INCORP. TG (ETIT)
$ws1.Range("A32").Value = 'ECOBANK TRANS. INCORP. TG (ETIT)'
$ws1.Range("B32").Value = 26
$ws1.Range("C32").Value = 20
$ws1.Range("D32").Value = 47.93
$ws1.Range("E32").Value = 5.88
$ws1.Range("G32").Value = '✅ Renforcer'

# Row 33: SOCIETE IVOIRIENNE DE BANQUE  (SIBC)
$ws1.Range("A33").Value = 'SOCIETE IVOIRIENNE DE BANQUE  (SIBC)'
$ws1.Range("B33").Value = 12
$ws1.Range("C33").Value = 3
$ws1.Range("D33").Value = 45.14
$ws1.Range("E33").Value = -3.36
$ws1.Range("G33").Value = '➖ Neutre'

# Row 35: BANK OF AFRICA SENEGAL (BOAS)
$ws1.Range("A35").Value = 'BANK OF AFRICA SENEGAL (BOAS)'
$ws1.Range("B35").Value = 10
$ws1.Range("C35").Value = 5
$ws1.Range("D35").Value = 39.73
$ws1.Range("E35").Value = 3.92
$ws1.Range("G35").Value = 'Non évalué'

# Row 36: SAPH CI (SPHC)
$ws1.Range("A36").Value = 'SAPH CI (SPHC)'
$ws1.Range("B36").Value = 13
$ws1.Range("C36").Value = 9
$ws1.Range("D36").Value = 36.42
$ws1.Range("E36").Value = 3.33

# Row 37: TOTALENERGIES MARKETING CI (TTLC)
$ws1.Range("A37").Value = 'TOTALENERGIES MARKETING CI (TTLC)'
$ws1.Range("B37").Value = 22
$ws1.Range("C37").Value = 20
$ws1.Range("D37").Value = 34.8
$ws1.Range("E37").Value = -3.7
$ws1.Range("G37").Value = '✅ Renforcer'

# Row 38: SOLIBRA CI (SLBC)
$ws1.Range("A38").Value = 'SOLIBRA CI (SLBC)'
$ws1.Range("B38").Value = 15
$ws1.Range("C38").Value = 11
$ws1.Range("D38").Value = 30.73
$ws1.Range("E38").Value = 7.3
$ws1.Range("G38").Value = '✅ Renforcer'

# Row 39: ECOBANK COTE D''IVOIRE (ECOC)
$ws1.Range("A39").Value = 'ECOBANK COTE D''''IVOIRE (ECOC)'
$ws1.Range("B39").Value = 8
$ws1.Range("C39").Value = 3
$ws1.Range("D39").Value = 27.8
$ws1.Range("E39").Value = 5
$ws1.Range("G39").Value = '➖ Neutre'

# Row 41: BANK OF AFRICA CI (BOAC)
$ws1.Range("G41").Value = 'Non évalué'

# Row 42: TOTALENERGIES MARKETING SN (TTLS)
$ws1.Range("A42").Value = 'TOTALENERGIES MARKETING SN (TTLS)'
$ws1.Range("D42").Value = 23.8
$ws1.Range("E42").Value = -3.41
$ws1.Range("G42").Value = 'Non évalué'

# Row 43: CFAO MOTORS CI (CFAC)
$ws1.Range("A43").Value = 'CFAO MOTORS CI (CFAC)'
$ws1.Range("C43").Value = 8
$ws1.Range("D43").Value = 21.06
$ws1.Range("E43").Value = -3.68
$ws1.Range("G43").Value = '➖ Neutre'

# Row 45: SAFCA CI (SAFC)
$ws1.Range("A45").Value = 'SAFCA CI (SAFC)'
$ws1.Range("B45").Value = 11
$ws1.Range("C45").Value = 7
$ws1.Range("D45").Value = 15.52
$ws1.Range("E45").Value = 6.47
$ws1.Range("G45").Value = '👀 À surveiller'

# Row 46: BANK OF AFRICA BN (BOAB)
$ws1.Range("A46").Value = 'BANK OF AFRICA BN (BOAB)'
$ws1.Range("B46").Value = 4
$ws1.Range("C46").Value = 3
$ws1.Range("D46").Value = 14.2
$ws1.Range("E46").Value = -1.9
$ws1.Range("G46").Value = 'Non évalué'

# Row 47: CIE CI (CIEC)
$ws1.Range("A47").Value = 'CIE CI (CIEC)'
$ws1.Range("B47").Value = 8
$ws1.Range("C47").Value = 7
$ws1.Range("D47").Value = 13.57
$ws1.Range("E47").Value = 5.75
$ws1.Range("F47").Value = '🟢 Achat'
$ws1.Range("G47").Value = '➖ Neutre'

# Row 48: AFRICA GLOBAL LOGISTICS CI (SDSC)
$ws1.Range("A48").Value = 'AFRICA GLOBAL LOGISTICS CI (SDSC)'
$ws1.Range("B48").Value = 2
$ws1.Range("C48").Value = 0
$ws1.Range("D48").Value = 10.33
$ws1.Range("E48").Value = 6.23
$ws1.Range("F48").Value = '🟡 Observer'
$ws1.Range("G48").Value = 'Non évalué'

# Row 51: VIVO ENERGY CI (SHEC)
$ws1.Range("C51").Value = 9
$ws1.Range("D51").Value = 7.41

# Row 55: PALM CI (PALC)
$ws1.Range("A55").Value = 'PALM CI (PALC)'
$ws1.Range("B55").Value = 8
$ws1.Range("C55").Value = 12
$ws1.Range("D55").Value = 0.4
$ws1.Range("E55").Value = 2.71
$ws1.Range("G55").Value = 'Non évalué'

# Row 56: TOTAL
$ws1.Range("A56").Value = 'TOTAL'
$ws1.Range("B56").Value = 0
$ws1.Range("C56").Value = 96
$ws1.Range("D56").Value = 0
$ws1.Range("E56").Value = 0
$ws1.Range("G56").Value = '➖ Neutre'

# Row 57: SICOR CI (SICC)
$ws1.Range("A57").Value = 'SICOR CI (SICC)'
$ws1.Range("B57").Value = 9
$ws1.Range("C57").Value = 11
$ws1.Range("D57").Value = -2.78
$ws1.Range("E57").Value = -4.35

# Row 58: SONATEL SN (SNTS)
$ws1.Range("C58").Value = 8
$ws1.Range("D58").Value = -6.8
$ws1.Range("G58").Value = '⚠️ Risque de décrochage'

# Row 59: LOTERIE NATIONALE DU BENIN (LNBB)
$ws1.Range("A59").Value = 'LOTERIE NATIONALE DU BENIN (LNBB)'
$ws1.Range("B59").Value = 0
$ws1.Range("C59").Value = 3
$ws1.Range("D59").Value = -7.03
$ws1.Range("E59").Value = -2.63
$ws1.Range("G59").Value = 'Non évalué'

# Row 60: NESTLE CI (NTLC)
$ws1.Range("A60").Value = 'NESTLE CI (NTLC)'
$ws1.Range("B60").Value = 1
$ws1.Range("C60").Value = 3
$ws1.Range("D60").Value = -10.22
$ws1.Range("E60").Value = 4.21

# Row 61: SODE CI (SDCC)
$ws1.Range("A61").Value = 'SODE CI (SDCC)'
$ws1.Range("B61").Value = 3
$ws1.Range("C61").Value = 7
$ws1.Range("D61").Value = -10.99
$ws1.Range("E61").Value = 6.56
$ws1.Range("G61").Value = '➖ Neutre'

# Row 62: BANK OF AFRICA BF (BOABF)
$ws1.Range("A62").Value = 'BANK OF AFRICA BF (BOABF)'
$ws1.Range("B62").Value = 9
$ws1.Range("C62").Value = 13
$ws1.Range("D62").Value = -12.42
$ws1.Range("E62").Value = -2.91
$ws1.Range("G62").Value = '➖ Neutre'

# Row 63: SOCIETE GENERALE COTE D'IVOIRE (SGBC)
$ws1.Range("A63").Value = 'SOCIETE GENERALE COTE D''IVOIRE (SGBC)'
$ws1.Range("B63").Value = 0
$ws1.Range("C63").Value = 9
$ws1.Range("D63").Value = -12.49
$ws1.Range("E63").Value = -1.14
$ws1.Range("G63").Value = 'Non évalué'

# Row 64: SUCRIVOIRE (SCRC)
$ws1.Range("A64").Value = 'SUCRIVOIRE (SCRC)'
$ws1.Range("B64").Value = 7
$ws1.Range("C64").Value = 16
$ws1.Range("D64").Value = -17.64
$ws1.Range("E64").Value = -1.5
$ws1.Range("G64").Value = '➖ Neutre'

# Row 65: SERVAIR ABIDJAN CI (ABJC)
$ws1.Range("A65").Value = 'SERVAIR ABIDJAN CI (ABJC)'
$ws1.Range("B65").Value = 7
$ws1.Range("C65").Value = 14
$ws1.Range("D65").Value = -23.32
$ws1.Range("E65").Value = 4.92
$ws1.Range("G65").Value = '✅ Renforcer'

# Row 66: AIR LIQUIDE CI (SIVC)
$ws1.Range("A66").Value = 'AIR LIQUIDE CI (SIVC)'
$ws1.Range("B66").Value = 3
$ws1.Range("C66").Value = 12
$ws1.Range("D66").Value = -27.61
$ws1.Range("E66").Value = -5.45
$ws1.Range("G66").Value = '➖ Neutre'

# Row 67: BANK OF AFRICA NG (BOAN)
$ws1.Range("A67").Value = 'BANK OF AFRICA NG (BOAN)'
$ws1.Range("B67").Value = 10
$ws1.Range("C67").Value = 17
$ws1.Range("D67").Value = -29.92
$ws1.Range("E67").Value = 4.17
$ws1.Range("G67").Value = 'Non évalué'

# Row 68: SICABLE CI (CABC)
$ws1.Range("A68").Value = 'SICABLE CI (CABC)'
$ws1.Range("B68").Value = 18
$ws1.Range("C68").Value = 30
$ws1.Range("D68").Value = -30.68
$ws1.Range("E68").Value = -7.5
$ws1.Range("G68").Value = '⚠️ Risque de décrochage'

# Row 69: SETAO CI (STAC)
$ws1.Range("C69").Value = 28
$ws1.Range("D69").Value = -33.87

# ----------------------------------------------------------------------------
# Sheet "Top_YTD"
# ----------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Top_YTD")

# Row 2
$ws2.Range("B2").Value = 209.27

# Row 3
$ws2.Range("B3").Value = 194.73

# Row 6
$ws2.Range("A6").Value = 'UNIWAX CI (UNXC)'
$ws2.Range("B6").Value = 56.55

# Row 7
$ws2.Range("A7").Value = 'SOCIETE IVOIRIENNE DE BANQUE  (SIBC)'
$ws2.Range("B7").Value = 54.39

# Row 8
$ws2.Range("A8").Value = 'BICI CI (BICC)'
$ws2.Range("B8").Value = 49.27

# Row 10
$ws2.Range("A10").Value = 'ECOBANK TRANS. INCORP. TG (ETIT)'
$ws2.Range("B10").Value = 46.56

Write-Output "BRVM data refresh applied."
